$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Cells in column D whose new values look like plain numbers. Excel would
# auto-convert a numeric-looking string assigned via .Value into a real
# number (losing the original text formatting, e.g. trailing zeros, and
# introducing floating point artifacts). These cells originally store the
# price as literal text, so force the number format to Text ("@") first
# so the assigned value is kept exactly as typed.
$textCells = @("D5","D6","D8","D9","D15","D16","D17","D18","D21","D22","D23","D24","D26","D27","D29","D31","D32","D34","D36","D41","D42","D43","D45","D46","D48","D49")
foreach ($addr in $textCells) {
    $ws.Range($addr).NumberFormat = "@"
}

# Apply the updated Price (column D) and Volume(1h) (column E) values.
$ws.Range('D2').Value = '67.879.90'
$ws.Range('E2').Value = '  -1.12%  '
$ws.Range('D3').Value = '2.430.34'
$ws.Range('E3').Value = '  -0.82%  '
$ws.Range('E4').Value = '  -0.04%  '
$ws.Range('D5').Value = '553.32'
$ws.Range('E5').Value = '  -0.51%  '
$ws.Range('D6').Value = '160.00'
$ws.Range('E6').Value = '  -0.59%  '
$ws.Range('E7').Value = '  +0.02%  '
$ws.Range('D8').Value = '0.509'
$ws.Range('E8').Value = '  +1.56%  '
$ws.Range('D9').Value = '0.159'
$ws.Range('E9').Value = '  +7.28%  '
$ws.Range('E10').Value = '  -0.66%  '
$ws.Range('E11').Value = '  -1.43%  '
$ws.Range('D13').Value = '67.783.59'
$ws.Range('E13').Value = '  -0.65%  '
$ws.Range('E14').Value = '  +1.37%  '
$ws.Range('D15').Value = '22.99'
$ws.Range('E15').Value = '  -1.24%  '
$ws.Range('D16').Value = '10.36'
$ws.Range('E16').Value = '  -3.24%  '
$ws.Range('D17').Value = '333.41'
$ws.Range('E17').Value = '  -1.96%  '
$ws.Range('D18').Value = '6.84'
$ws.Range('E18').Value = '  -2.04%  '
$ws.Range('E20').Value = '  -0.25%  '
$ws.Range('D21').Value = '1.85'
$ws.Range('E21').Value = '  +0.35%  '
$ws.Range('D22').Value = '66.11'
$ws.Range('E22').Value = '  -0.43%  '
$ws.Range('D23').Value = '3.62'
$ws.Range('E23').Value = '  -0.23%  '
$ws.Range('D24').Value = '8.10'
$ws.Range('E24').Value = '  +0.40%  '
$ws.Range('D25').Value = '0.0₃0807'
$ws.Range('E25').Value = '  -0.54%  '
$ws.Range('D26').Value = '7.08'
$ws.Range('E26').Value = '  -0.57%  '
$ws.Range('D27').Value = '0.999'
$ws.Range('E27').Value = '  -0.09%  '
$ws.Range('E28').Value = '  +0.55%  '
$ws.Range('D29').Value = '416.08'
$ws.Range('E29').Value = '  -4.16%  '
$ws.Range('E30').Value = '  -0.36%  '
$ws.Range('D31').Value = '160.39'
$ws.Range('E31').Value = '  +2.72%  '
$ws.Range('D32').Value = '18.94'
$ws.Range('E32').Value = '  -0.56%  '
$ws.Range('E33').Value = '  -0.09%  '
$ws.Range('D34').Value = '17.83'
$ws.Range('E35').Value = '  -3.47%  '
$ws.Range('D36').Value = '0.294'
$ws.Range('E36').Value = '  -2.30%  '
$ws.Range('E37').Value = '  -3.12%  '
$ws.Range('E38').Value = '  +0.69%  '
$ws.Range('E39').Value = '  -1.03%  '
$ws.Range('E40').Value = '  -1.92%  '
$ws.Range('D41').Value = '3.32'
$ws.Range('E41').Value = '  +0.03%  '
$ws.Range('D42').Value = '129.51'
$ws.Range('E42').Value = '  -1.69%  '
$ws.Range('D43').Value = '0.0705'
$ws.Range('E43').Value = '  -1.06%  '
$ws.Range('E44').Value = '  -0.41%  '
$ws.Range('D45').Value = '0.555'
$ws.Range('E45').Value = '  -0.60%  '
$ws.Range('D46').Value = '0.0913'
$ws.Range('E46').Value = '  +0.81%  '
$ws.Range('E47').Value = '  +0.39%  '
$ws.Range('D48').Value = '1.34'
$ws.Range('E48').Value = '  -6.53%  '
$ws.Range('D49').Value = '16.51'
$ws.Range('E49').Value = '  -1.54%  '
$ws.Range('E50').Value = '  +2.99%  '
